# Fruta / hortaliza, semanal
# Insert two new weekly data rows (Mango, Agrícola del Norte S.A. de Arica)
# right before the existing row 47, shifting all following rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 47 (pushes old row 47 -> 49, old row 93 -> 95)
$ws.Rows.Item(47).Insert()
$ws.Rows.Item(47).Insert()

# --- New row 47 ---
$ws.Range("A47").Value2 = 1
$ws.Range("B47").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C47").Value2 = "Arica y Parinacota"
$ws.Range("D47").Value2 = 44482
$ws.Range("E47").Value2 = 15
$ws.Range("F47").Value2 = "Fruta"
$ws.Range("G47").Value2 = 100108
$ws.Range("H47").Value2 = "Tropicales y subtropicales"
$ws.Range("I47").Value2 = 100108002
$ws.Range("J47").Value2 = "Mango"
$ws.Range("K47").Value2 = "Sin especificar"
$ws.Range("L47").Value2 = "Especial"
$ws.Range("M47").Value2 = 456
$ws.Range("N47").Value2 = 4500
$ws.Range("O47").Value2 = 5000
$ws.Range("P47").Value2 = 4750
$ws.Range("Q47").Value2 = "$/bandeja 4 kilos"
$ws.Range("R47").Value2 = "Perú"
$ws.Range("S47").Value2 = 1188
$ws.Range("T47").Value2 = 4

# --- New row 48 ---
$ws.Range("A48").Value2 = 1
$ws.Range("B48").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C48").Value2 = "Arica y Parinacota"
$ws.Range("D48").Value2 = 44482
$ws.Range("E48").Value2 = 15
$ws.Range("F48").Value2 = "Fruta"
$ws.Range("G48").Value2 = 100108
$ws.Range("H48").Value2 = "Tropicales y subtropicales"
$ws.Range("I48").Value2 = 100108002
$ws.Range("J48").Value2 = "Mango"
$ws.Range("K48").Value2 = "Sin especificar"
$ws.Range("L48").Value2 = "Primera"
$ws.Range("M48").Value2 = 456
$ws.Range("N48").Value2 = 4500
$ws.Range("O48").Value2 = 5000
$ws.Range("P48").Value2 = 4750
$ws.Range("Q48").Value2 = "$/bandeja 4 kilos"
$ws.Range("R48").Value2 = "Perú"
$ws.Range("S48").Value2 = 1188
$ws.Range("T48").Value2 = 4

# Make sure the date cells keep the same date number format as the rest of column D
$ws.Range("D47").NumberFormat = $ws.Range("D46").NumberFormat
$ws.Range("D48").NumberFormat = $ws.Range("D46").NumberFormat
